$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the item code / name fields (shared strings SB008... -> SB009...)
$ws.Range("A2").Value = "SB009-SECOND"
$ws.Range("B2").Value = "SB009"
$ws.Range("P2").Value = "SB009-SECOND"
$ws.Range("T2").Value = "SB009"

# Update the unit sale price
$ws.Range("G2").Value = 55
